# Commit: [base] - [`outputToCloud(resource)`]: support the transferring of
# output artifact to the cloud.
#
# This adds a new "base" category function `outputToCloud(resource)` and a
# new "text" target category (function list) whose sole member is
# `spellCheck(var,profile,text)`, on the hidden "#system" sheet that backs
# the workbook's data-validation / autocomplete defined-name lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Insert a brand-new column at "Y" (shifts old Y..AD -> Z..AE, i.e. the
#    "web", "webalert", "webcookie", "ws", "ws.async", "xml" lists each move
#    one column to the right) and seed it with the new "text" target list
#    (header + single function name).
# ---------------------------------------------------------------------------
$ws.Columns("Y:Y").Insert()
$ws.Cells.Item(1, 25).Value2 = "text"
$ws.Cells.Item(2, 25).Value2 = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------------
# 2) Insert "outputToCloud(resource)" into the "base" function list (column
#    E), keeping it alphabetically sorted: it lands on row 22, between
#    "macro(file,sheet,name)" (row 21) and "prependText(var,prependWith)"
#    (formerly row 22, now row 23). Only column E shifts -- the other
#    category columns on the same rows are untouched.
# ---------------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $moved = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r + 1, 5).Value2 = $moved
}
$ws.Cells.Item(22, 5).Value2 = "outputToCloud(resource)"

# ---------------------------------------------------------------------------
# 3) Insert the new "text" category name into the "target" list (column A),
#    keeping it alphabetically sorted: it lands on row 25, between "step"
#    (row 24) and "web" (formerly row 25, now row 26).
# ---------------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $moved = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 1).Value2 = $moved
}
$ws.Cells.Item(25, 1).Value2 = "text"

# ---------------------------------------------------------------------------
# 4) Fix up the workbook-level defined names so they keep pointing at the
#    right ranges after the shifts above, and register the new "text" name.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = '=''#system''!$E$2:$E$39'
$wb.Names.Item("target").RefersTo = '=''#system''!$A$2:$A$31'
$wb.Names.Item("web").RefersTo = '=''#system''!$Z$2:$Z$129'
$wb.Names.Item("webalert").RefersTo = '=''#system''!$AA$2:$AA$8'
$wb.Names.Item("webcookie").RefersTo = '=''#system''!$AB$2:$AB$8'
$wb.Names.Item("ws").RefersTo = '=''#system''!$AC$2:$AC$17'
$wb.Names.Item("ws.async").RefersTo = '=''#system''!$AD$2:$AD$8'
$wb.Names.Item("xml").RefersTo = '=''#system''!$AE$2:$AE$27'
$wb.Names.Add("text", '=''#system''!$Y$2:$Y$2')
